$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New row 5: 1148. Article Views I
# ---------------------------------------------------------------------------
$ws.Range("A5").Value2 = "1148. Article Views I"
$ws.Range("B5").Value2 = $ws.Range("B4").Value2
$ws.Range("C5").Value2 = $ws.Range("C4").Value2
$ws.Range("D5").Value2 = "SELECT DISTINCT col AS id FROM View WHERE id1 = id2 ORDER BY id ASC;"
$ws.Range("E5").Value2 = "https://leetcode.com/problems/article-views-i/solutions/1945221/oracle-sql-simple-query-using-where-clause/?envType=study-plan-v2&envId=top-sql-50 "

# ---------------------------------------------------------------------------
# New row 6: 1683. Invalid Tweets
# ---------------------------------------------------------------------------
$ws.Range("A6").Value2 = "1683. Invalid Tweets"
$ws.Range("B6").Value2 = $ws.Range("B4").Value2
$ws.Range("C6").Value2 = $ws.Range("C4").Value2
$ws.Range("D6").Value2 = "Use LENGTH(col) > val"
$ws.Range("E6").Value2 = "https://leetcode.com/problems/invalid-tweets/solutions/3857245/100-easy-fast-clean-one-line-solution/?envType=study-plan-v2&envId=top-sql-50 "

# Match the formatting (styles) of the preceding data row for the two new rows
$ws.Range("A4:E4").Copy()
$ws.Range("A5:E5").PasteSpecial(-4122)
$ws.Range("A6:E6").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Wire up the Link column hyperlinks (external relationships)
$ws.Hyperlinks.Add($ws.Range("E5"), "https://leetcode.com/problems/article-views-i/solutions/1945221/oracle-sql-simple-query-using-where-clause/?envType=study-plan-v2&envId=top-sql-50 ")
$ws.Hyperlinks.Add($ws.Range("E6"), "https://leetcode.com/problems/invalid-tweets/solutions/3857245/100-easy-fast-clean-one-line-solution/?envType=study-plan-v2&envId=top-sql-50 ")

# Adding a hyperlink re-applies formatting - restore the row formats again
$ws.Range("A4:E4").Copy()
$ws.Range("A5:E5").PasteSpecial(-4122)
$ws.Range("A6:E6").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# Grow Table2 to cover the two new rows
# ---------------------------------------------------------------------------
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:E6"))

# ---------------------------------------------------------------------------
# Match the author's final cursor/view state
# ---------------------------------------------------------------------------
$ws.Range("D15").Select()
